# Bugfix in copyAndInsert. copyCellsFrom(srcRow) supports now merged regions
# (inside the current row).
#
# This reproduces, via Excel COM automation, the two example edits that
# exercise the fix:
#   - sheet "sheet 1": row 4 gets a new merged "combined" cell (F4:H4),
#     mirroring the existing merged "combined" cell structure used on the
#     "Row-actions" sheet.
#   - sheet "Row-actions": a brand-new row (row 2) is inserted above the
#     existing sample rows, itself containing a merged cell (H2:I2), to show
#     that copying a source row that contains a merge now works.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# sheet 1 : add merged "combined" cell at F4:H4
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("sheet 1")

$ws1.Range("F4").Value = "combined"
$ws1.Range("F4:H4").HorizontalAlignment = -4131
$ws1.Range("F4:H4").Merge()

# ---------------------------------------------------------------------
# Row-actions : insert a new row of sample data (with its own merged
# cell H2:I2) right after the header/first row.
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Row-actions")

# Copy the formatting of the existing sample row (row 3) onto row 2 first
# (this is the "copyCellsFrom(srcRow)" codepath the bugfix targets), then
# overwrite the values for the new row.
$ws5.Range("C3:J3").Copy()
$ws5.Range("C2").PasteSpecial(-4122)
$ws5.Range("A2:J2").RowHeight = $ws5.Range("A3:J3").RowHeight

$ws5.Range("C2").Value = "Test 1"
$ws5.Range("D2").Value = 0.01
$ws5.Range("E2").Value = 15
$ws5.Range("F2").Value = "Box 1"
$ws5.Range("G2").Value = "red, underline 1"
$ws5.Range("H2").Value = "combined 1"
$ws5.Range("J2").Value = "last 1"

$ws5.Range("H2:I2").Merge()

# ---------------------------------------------------------------------
# Restore "sheet 1" as the active sheet/tab (it was "Row-actions" before).
# ---------------------------------------------------------------------
$ws1.Activate()
